$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 3 (Short Att, Short Comp, Deep Att, Deep Comp, Short Int)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 189
$wsOff.Range("C3").Value = 138
$wsOff.Range("D3").Value = 53
$wsOff.Range("E3").Value = 30
$wsOff.Range("F3").Value = 5

# Sheet "DEF" - row 3 (Short Att, Short Comp, Deep Att, Deep Comp, Short Int)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 148
$wsDef.Range("C3").Value = 104
$wsDef.Range("D3").Value = 44
$wsDef.Range("E3").Value = 25
$wsDef.Range("F3").Value = 3
